# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2020-09-30 00:00:00"

$ws.Range("O2").Value = 114273999.06
$ws.Range("P2").Value = 288355225.66
$ws.Range("Q2").Value = 174759025.36
$ws.Range("R2").Value = 49.8011535162
$ws.Range("S2").Value = 107843086.31
$ws.Range("T2").Value = 107843086.31
$ws.Range("U2").Value = 53.243851916
$ws.Range("V2").Value = 13048531.3
$ws.Range("W2").Value = 25321119.14
$ws.Range("X2").Value = 1073958.91
$ws.Range("Y2").Value = 135112986.64
$ws.Range("Z2").Value = 134375750.53
$ws.Range("AA2").Value = 20101751.47
$ws.Range("AG2").Value = 2236812.94
$ws.Range("AP2").Value = 51.9361302775
$ws.Range("AQ2").Value = 63.645843632546
$ws.Range("AR2").Value = 61.79
$ws.Range("AS2").Value = 95628918.17
$ws.Range("AT2").Value = 62.16418526909
